# remove calculation of Bierkruglauf
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PktBierkruglauf")

# The sheet no longer calculates StrafPkt (C) or Zeit (D) - clear their
# data but keep the column formatting in place.
$ws.Range("C1:D4").ClearContents() | Out-Null

# Rename the remaining header from "Pkt" to "Gesamt" (it's now the only /
# overall total column).
$ws.Range("B1").Value = "Gesamt"

# PktBierkruglauf becomes the active/selected sheet (previously PktWackelturm).
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
